# Weekly price-data update: insert the latest week's two new records
# (quality grades "Primera" and "Segunda") for Crespo record cabbage at
# Terminal Hortofrutícola Agro Chillán, pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 265-266; everything from the old row 265
# onward shifts down by two rows (old 365/366 become new 367/368).
$ws.Rows("265:266").Insert()

# --- New row 265 : Primera ---
$ws.Range("A265").Value = 7
$ws.Range("B265").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C265").Value = "Ñuble"
$ws.Range("D265").Value = 45027
$ws.Range("E265").Value = 16
$ws.Range("F265").Value = 100112006
$ws.Range("G265").Value = "Repollo"
$ws.Range("H265").Value = "Crespo record"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 300
$ws.Range("K265").Value = 1300
$ws.Range("L265").Value = 1400
$ws.Range("M265").Value = 1350
$ws.Range("N265").Value = "`$/unidad"
$ws.Range("O265").Value = "Provincia de Diguillín"
$ws.Range("P265").Value = 1350
$ws.Range("Q265").Value = 1
$ws.Range("R265").Value = "Hortaliza"

# --- New row 266 : Segunda ---
$ws.Range("A266").Value = 7
$ws.Range("B266").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C266").Value = "Ñuble"
$ws.Range("D266").Value = 45027
$ws.Range("E266").Value = 16
$ws.Range("F266").Value = 100112006
$ws.Range("G266").Value = "Repollo"
$ws.Range("H266").Value = "Crespo record"
$ws.Range("I266").Value = "Segunda"
$ws.Range("J266").Value = 300
$ws.Range("K266").Value = 1000
$ws.Range("L266").Value = 1000
$ws.Range("M266").Value = 1000
$ws.Range("N266").Value = "`$/unidad"
$ws.Range("O266").Value = "Provincia de Diguillín"
$ws.Range("P266").Value = 1000
$ws.Range("Q266").Value = 1
$ws.Range("R266").Value = "Hortaliza"
